$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.531.30"
$ws.Range("E2").Value = "  -6.48%  "
$ws.Range("D3").Value = "3.314.88"
$ws.Range("E3").Value = "  -3.28%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'561.42"
$ws.Range("E5").Value = "  -3.52%  "
$ws.Range("D6").Value = "'129.33"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "3.320.14"
$ws.Range("E8").Value = "  -3.05%  "
$ws.Range("D9").Value = "'0.470"
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("D10").Value = "'7.39"
$ws.Range("E10").Value = "  -3.08%  "
$ws.Range("D11").Value = "'0.117"
$ws.Range("E11").Value = "  -6.03%  "
$ws.Range("D12").Value = "'0.372"
$ws.Range("E12").Value = "  -3.09%  "
$ws.Range("D13").Value = "3.877.76"
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "3.309.53"
$ws.Range("E15").Value = "  -3.63%  "
$ws.Range("D16").Value = "'0.0000167"
$ws.Range("E16").Value = "  -5.83%  "
$ws.Range("D17").Value = "'24.37"
$ws.Range("E17").Value = "  -2.47%  "
$ws.Range("D18").Value = "59.764.96"
$ws.Range("E18").Value = "  -6.03%  "
$ws.Range("D19").Value = "'5.62"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").Value = "'13.29"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "'8.97"
$ws.Range("E21").Value = "  -9.26%  "
$ws.Range("D22").Value = "'351.45"
$ws.Range("E22").Value = "  -8.94%  "
$ws.Range("D23").Value = "'0.553"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").Value = "3.441.13"
$ws.Range("E25").Value = "  -3.56%  "
$ws.Range("D26").Value = "'68.58"
$ws.Range("E26").Value = "  -7.15%  "
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D29").Value = "'7.38"
$ws.Range("E29").Value = "  +4.83%  "
$ws.Range("D30").Value = "'1.48"
$ws.Range("E30").Value = "  +4.24%  "
$ws.Range("D31").Value = "'7.82"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").Value = "'0.152"
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("D33").Value = "'2.10"
$ws.Range("E33").Value = "  -5.08%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "3.339.56"
$ws.Range("E35").Value = "  -3.41%  "
$ws.Range("D36").Value = "'22.73"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "'5.32"
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("D38").Value = "'6.78"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'1.49"
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("D40").Value = "'158.01"
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("D41").Value = "'0.0750"
$ws.Range("E41").Value = "  -3.23%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("B43").Value = "ONDO"
$ws.Range("C43").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D43").Value = "'1.19"
$ws.Range("E43").Value = "  +6.93%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'40.72"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").Value = "'4.32"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").Value = "'0.743"
$ws.Range("E46").Value = "  -5.62%  "
$ws.Range("D47").Value = "'23.27"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "'1.55"
$ws.Range("E48").Value = "  -4.01%  "
$ws.Range("D49").Value = "'6.73"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'22.15"
$ws.Range("E50").Value = "  +8.83%  "
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'2.41"
$ws.Range("E51").Value = "  +15.67%  "
